# Edit script: "Celestial Symphony" (astronomy) essay -> "Journeying Through
# the Realm of Chemistry" essay, plus author/email swap and a font-name fix
# (TimesNewToman -> Times New Roman) across every run in the document.

$d = $word.ActiveDocument

$ok = $d.Content.Find.Execute("The Celestial Symphony: Unveiling the Harmony of the Cosmos", $true, $false, $false, $false, $false, $true, 1, $false, "Journeying Through the Realm of Chemistry: Unveiling the Secrets of Matter", 2)
if (-not $ok) { throw "replace 0 failed: The Celestial Symphony: Unveiling the Ha" }
$ok = $d.Content.Find.Execute("Amelia Carter", $true, $false, $false, $false, $false, $true, 1, $false, "Emily Harper", 2)
if (-not $ok) { throw "replace 1 failed: Amelia Carter" }
$ok = $d.Content.Find.Execute("ameliacarter@spaceobservatory", $true, $false, $false, $false, $false, $true, 1, $false, "emilyharper0622@ymail", 2)
if (-not $ok) { throw "replace 2 failed: ameliacarter@spaceobservatory" }
$ok = $d.Content.Find.Execute("edu", $true, $true, $false, $false, $false, $true, 1, $false, "net", 2)
if (-not $ok) { throw "replace 3 failed: edu" }
$ok = $d.Content.Find.Execute("In the vast expanse of the cosmos, celestial bodies dance to an intricate rhythm, creating a symphony of cosmic wonders", $true, $false, $false, $false, $false, $true, 1, $false, "In the vast expanse of human knowledge, chemistry stands as a captivating language that unlocks the hidden secrets of matter", 2)
if (-not $ok) { throw "replace 4 failed: In the vast expanse of the cosmos, celes" }
$ok = $d.Content.Find.Execute(" The universe, like a celestial orchestra, weaves together the melodies of stars, planets, and galaxies, inviting us to decipher their harmonies and unravel the mysteries of existence", $true, $false, $false, $false, $false, $true, 1, $false, " It is a science that investigates the fundamental building blocks of the universe and the intricate interactions between them, guiding us towards an understanding of the world around us", 2)
if (-not $ok) { throw "replace 5 failed:  The universe, like a celestial orchestr" }
$ok = $d.Content.Find.Execute(" From the blazing fires of the sun to the gentle glow of distant nebulas, each cosmic entity contributes to the grand composition, echoing the interconnectedness of all things", $true, $false, $false, $false, $false, $true, 1, $false, " From the smallest atoms to the sprawling galaxies, chemistry holds the key to unraveling the mysteries of existence", 2)
if (-not $ok) { throw "replace 6 failed:  From the blazing fires of the sun to th" }
$ok = $d.Content.Find.Execute("As we embark on this odyssey of cosmic exploration, we are granted a glimpse into the profound beauty and complexity of the universe", $true, $false, $false, $false, $false, $true, 1, $false, "Delving into the realm of chemistry, we embark on an extraordinary quest to explore the nature of substances and their transformations", 2)
if (-not $ok) { throw "replace 7 failed: As we embark on this odyssey of cosmic e" }
$ok = $d.Content.Find.Execute(" We witness the birth and death of stars, the graceful ballet of planets around their suns, and the explosive spectacle of supernovae that herald the creation of new elements", $true, $false, $false, $false, $false, $true, 1, $false, " We learn about the elements, the basic units of matter, and the ways in which they combine to form compounds with diverse properties", 2)
if (-not $ok) { throw "replace 8 failed:  We witness the birth and death of stars" }
$ok = $d.Content.Find.Execute(" With each observation, we deepen our understanding of the universe's composition, its evolution, and our place within this cosmic tapestry", $true, $false, $false, $false, $false, $true, 1, $false, " Through chemical reactions, we witness the dynamic interplay of atoms and molecules, as they rearrange and recombine, creating new substances with distinct characteristics", 2)
if (-not $ok) { throw "replace 9 failed:  With each observation, we deepen our un" }
$ok = $d.Content.Find.Execute("The study of celestial phenomena has ignited human curiosity for millennia, inspiring profound contemplations about our origins, our destiny, and the nature of reality itself", $true, $false, $false, $false, $false, $true, 1, $false, "Unraveling the enigmas of chemistry enables us to decipher the fundamental principles that govern the behavior of matter", 2)
if (-not $ok) { throw "replace 10 failed: The study of celestial phenomena has ign" }
$ok = $d.Content.Find.Execute(" From ancient astronomers who charted the movements of celestial bodies to modern astrophysicists who probe the depths of space with cutting-edge telescopes, humanity's quest for knowledge about the cosmos has been an enduring endeavor, revealing the intricate workings of the universe and expanding our horizons of understanding", $true, $false, $false, $false, $false, $true, 1, $false, " We uncover the laws of thermodynamics, which dictate the flow of energy and the direction of chemical change", 2)
if (-not $ok) { throw "replace 11 failed:  From ancient astronomers who charted th" }
$ok = $d.Content.Find.Execute("The celestial symphony that unfolds before us is a testament to the interconnectedness of the universe and the profound beauty that exists beyond our terrestrial sphere", $true, $false, $false, $false, $false, $true, 1, $false, "In this exploration of chemistry, we have embarked on a captivating journey through the realm of matter, unveiling the secrets of its structure, properties, and transformations", 2)
if (-not $ok) { throw "replace 12 failed: The celestial symphony that unfolds befo" }
$ok = $d.Content.Find.Execute(" As we continue to unravel the mysteries of the cosmos, we gain a deeper appreciation for our place within the vastness of existence and the intricate harmonies that govern the symphony of the stars", $true, $false, $false, $false, $false, $true, 1, $false, " We have delved into the fundamental principles that govern chemical reactions and the intricate interactions between atoms and molecules", 2)
if (-not $ok) { throw "replace 13 failed:  As we continue to unravel the mysteries" }
$ok = $d.Content.Find.Execute(" The study of celestial phenomena ignites our imagination, inspiring awe and wonder at the boundless marvels of the universe, reminding us that we are part of a grand cosmic dance that has been playing out for eons and will continue long after our own brief existence", $true, $false, $false, $false, $false, $true, 1, $false, " Through this exploration, we have gained a deeper understanding of the world around us and the remarkable complexity of the universe we inhabit", 2)
if (-not $ok) { throw "replace 14 failed:  The study of celestial phenomena ignite" }

# --- Insert the three new sentences that were added to the first body
# paragraph, right after the "...direction of chemical change" sentence
# (which itself replaced the old "...horizons of understanding" sentence
# above) and before the paragraph's closing period.
$tail = $d.Content
$found = $tail.Find.Execute(" We uncover the laws of thermodynamics, which dictate the flow of energy and the direction of chemical change", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not locate the thermodynamics sentence to append after" }
$tail.Collapse(0)
$tail.InsertAfter(".")
$tail.Collapse(0)
$tail.InsertAfter(" We unravel the mysteries of chemical bonding, the forces that hold atoms together and determine the properties of compounds")
$tail.Collapse(0)
$tail.InsertAfter(".")
$tail.Collapse(0)
$tail.InsertAfter(" Moreover, we delve into the intricacies of chemical reactions, exploring the factors that influence their rates and the mechanisms by which they occur")

# --- Fix the misspelled font name "TimesNewToman" -> "Times New Roman" on
# every run in the document. Applying Font.Name through a freshly bounded
# Range per paragraph (rather than through Paragraph.Range directly) keeps
# the paragraph-mark run properties (w:pPr/w:rPr) untouched.
foreach ($p in $d.Paragraphs) {
    $s = $p.Range.Start
    $e = $p.Range.End
    if ($e -gt $s) {
        $rr = $d.Range($s, $e)
        $rr.Font.Name = "Times New Roman"
    }
}

# --- Add a new empty paragraph right after the closing "Summary" paragraph,
# just before the section break. Doing this after the font fix above means
# the fresh paragraph mark naturally inherits "Times New Roman".
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
